$wb = $excel.ActiveWorkbook

# --- Rename header cells in existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$data = @(
    @(44934.99999999999, 10, 9.999999987584344, 10.00000001332509),
    @(44941.99999999999, 10, 9.999999986475075, 10.00000001295745),
    @(44948.99999999999, 10, 9.999999918168481, 10.00000008752072),
    @(44955.99999999999, 10, 9.999999757757761, 10.00000026274956),
    @(44962.99999999999, 10, 9.999999520086535, 10.00000053206071),
    @(44969.99999999999, 10, 9.999999253548637, 10.00000083570027),
    @(44976.99999999999, 10, 9.999998931190307, 10.00000118060938),
    @(44983.99999999999, 10, 9.999998502097837, 10.00000154220779),
    @(44990.99999999999, 10, 9.999998146947442, 10.00000199539892),
    @(44997.99999999999, 10, 9.999997762431832, 10.00000242710333)
)

$row = 2
foreach ($rec in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $rec[0]
    $wsForecast.Cells.Item($row, 2).Value = $rec[1]
    $wsForecast.Cells.Item($row, 3).Value = $rec[2]
    $wsForecast.Cells.Item($row, 4).Value = $rec[3]
    $row++
}

# --- Match formatting to the sibling sheets: bold/centered header style,
#     and the date number-format style used for the "ds" date column ---
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
